$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Date Created" column header to "Date Issued"
$ws.Range("C2").Value = "Date Issued"

# Move the current selection to C3, matching the validated cell
$ws.Range("C3").Select()
